$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking "Price" (D) and "Volume(1h)" (E) columns are plain text cells.
# A handful of Price values look like plain numbers ("198.76", "0.202", ...);
# a leading apostrophe forces Excel to keep them as text (matching the
# original inlineStr cells) instead of silently parsing them into numbers.

$ws.Range("D2").Value = "76.434.02"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "2.958.14"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'198.76"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").Value = "'597.11"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "'0.202"
$ws.Range("E9").Value = "  +5.49%  "
$ws.Range("D10").Value = "2.958.16"
$ws.Range("E10").Value = "  +2.72%  "
$ws.Range("D11").Value = "'0.445"
$ws.Range("E11").Value = "  +12.43%  "
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("E13").Value = "  +4.48%  "
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").Value = "'28.54"
$ws.Range("E15").Value = "  +4.74%  "
$ws.Range("D16").Value = "76.387.57"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "2.952.97"
$ws.Range("E18").Value = "  +2.51%  "
$ws.Range("D19").Value = "'13.71"
$ws.Range("E19").Value = "  +9.37%  "
$ws.Range("D20").Value = "'8.77"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "'377.20"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  +4.63%  "
$ws.Range("D24").Value = "'72.68"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").Value = "3.093.34"
$ws.Range("E26").Value = "  +2.13%  "
$ws.Range("D27").Value = "'4.32"
$ws.Range("E27").Value = "  +2.35%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("E31").Value = "  +10.63%  "
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("D33").Value = "'496.57"
$ws.Range("E33").Value = "  -2.16%  "
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'165.68"
$ws.Range("E36").Value = "  +1.52%  "
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("E38").Value = "  +14.29%  "
$ws.Range("E39").Value = "  +18.64%  "
$ws.Range("D40").Value = "'19.97"
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("D41").Value = "'0.112"
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "'180.61"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").Value = "'4.94"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("D46").Value = "'39.86"
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("D48").Value = "'0.592"
$ws.Range("E48").Value = "  +2.40%  "
$ws.Range("E49").Value = "  +4.67%  "
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("E51").Value = "  +0.73%  "
